# Applies the OOXML diff: Persian text fixes (spacing, hamza, run-split
# normalization) plus a numbering-format change (decimal -> lowerLetter)
# for list level 4 of the list used by paragraphs 31-35.

$d = $word.ActiveDocument

# --- Text corrections (ZWNJ/half-space fixes, "tayid" -> "ta'yid", etc.) ---
$replacements = @(
    @("با واردکردن نام، نام‌خانوادگی، نام کاربری، آدرس ایمیل، رمز عبور و تکرار آن، گزینه ثبت‌نام برای کاربر جدید از رنگ قرمز به رنگ سبز درمی‌آید که نشان‌دهنده امکان ثبت‌نام اولیه برای اوست", "با واردکردن نام، نام خانوادگی، نام کاربری، آدرس ایمیل، رمز عبور و تکرار آن، گزینه ثبت‌نام برای کاربر جدید از رنگ قرمز به رنگ سبز درمی‌آید که نشان‌دهنده امکان ثبت‌نام اولیه برای اوست"),
    @("پس از وارد کردن اطلاعات و صحت‌سنجی آن‌ها در مرورگر(شماره تلفن نامعتبر، نام نامعتبر، SQL Injection) این اطلاعات بک‌اند فرستاده می‌شود.", "پس از واردکردن اطلاعات و صحت‌سنجی آن‌ها در مرورگر (شماره‌تلفن نامعتبر، نام نامعتبر، SQL Injection) این اطلاعات بک‌اند فرستاده می‌شود."),
    @("کاربر با مراجعه به ایمیل خود، ایمیل حاوی لینک تایید حساب کاربری خود را باز و به صفحه مربوطه هدایت می‌شود.", "کاربر با مراجعه به ایمیل خود، ایمیل حاوی لینک تأیید حساب کاربری خود را باز و به صفحه مربوطه هدایت می‌شود."),
    @("پس از بازدید سایت مخصوص به کاربر پیغام تایید به بک‌اند فرستاده می‌شود و حساب کاربری به حالت فعال در‌می‌آید و این موضوع در پیشخوان کاربر قابل مشاهده می‌شود.", "پس از بازدید سایت مخصوص به کاربر پیغام تأیید به بک‌اند فرستاده می‌شود و حساب کاربری به حالت فعال درمی‌آید و این موضوع در پیشخوان کاربر قابل‌مشاهده می‌شود."),
    @("در پنجره بازشده مبنی بر تایید خرید اشتراک، آن را تایید یا رد می‌کند.", "در پنجره بازشده مبنی بر تأیید خرید اشتراک، آن را تأیید یا رد می‌کند."),
    @("در صورت تایید، به درگاه پرداخت منتقل می‌شود.", "در صورت تأیید، به درگاه پرداخت منتقل می‌شود."),
    @("مجددا وارد صفحه ورود می‌شود و با نام کاربری و رمز عبور جدید وارد میزکار خود می‌شود.", "مجدداً وارد صفحه ورود می‌شود و با نام کاربری و رمز عبور جدید وارد میز کار خود می‌شود."),
    @("اگر کاربر رمز عبور را فراموش کرده باشد وارد آن قسمت می‌شود که به تفصیل توضیح داده شده است.", "اگر کاربر رمز عبور را فراموش کرده باشد وارد آن قسمت می‌شود که به‌تفصیل توضیح داده شده است."),
    @("اطلاعات موردنظر خود را تغییر یا از این اقدام صرف نظر می‌کند.", "اطلاعات موردنظر خود را تغییر یا از این اقدام صرف‌نظر می‌کند."),
    @("خلاصه‌ای از تغییرات به او نمایش داده می‌شود تا کاربر آن را تایید کند.", "خلاصه‌ای از تغییرات به او نمایش داده می‌شود تا کاربر آن را تأیید کند."),
    @("پس از تایید، پیغامی مبنی بر موفقیت‌آمیز بودن این عمل و همچنین ایمیلی با همین عنوان به اطلاع او می‌رسد.", "پس از تأیید، پیغامی مبنی بر موفقیت‌آمیز بودن این عمل و همچنین ایمیلی با همین عنوان به اطلاع او می‌رسد."),
    @("در این صفحه، کاربر به مشخصات آگهی‌ که مالک آن، ‌یکی دیگر از کاربران است، دسترسی دارد.", "در این صفحه، کاربر به مشخصات آگهی که مالک آن، یکی دیگر از کاربران است، دسترسی دارد."),
    @("دیگر مشخصات خودرو که به صورت اسلاید بار، گزینه، یا جعبه متن هست را مشخص می‌کند.", "دیگر مشخصات خودرو که به‌صورت اسلاید بار، گزینه، یا جعبه متن هست را مشخص می‌کند."),
    @("مشخات جزئی خودرو را در صورت تمایل تغییر می‌دهد.", "مشخصات جزئی خودرو را در صورت تمایل تغییر می‌دهد."),
    @("در هر مرحله لیستی از آگهی‌های متناسب با فیلترها به کاربر نمایش داده‌می‌شود و کاربر می‌تواند فیلترهارا اضافه یا کم و یا ویرایش کند.", "در هر مرحله لیستی از آگهی‌های متناسب با فیلترها به کاربر نمایش داده می‌شود و کاربر می‌تواند فیلترها را اضافه یا کم و یا ویرایش کند.")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found -> $old"
    }
}

# --- Numbering: abstractNumId 7 / numId 23, level 4 (w:ilvl=3) ---
# decimal ("1.", "2.", ...) -> lowerLetter ("a.", "b.", ...)
$lp = $d.Paragraphs(31)
$lvl4 = $lp.Range.ListFormat.ListTemplate.ListLevels.Item(4)
$lvl4.NumberStyle = 4
$lvl4.NumberFormat = "%4."

Write-Output "done"